$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D. Excel shifts the existing D:K data to
# E:L (carrying over each cell's style), and the newly created column D
# cells default to the general style.
$ws.Columns("D").Insert()

# Populate the new column D with the FY2018 figures, one row at a time.
# Each cell's number format is first copied over from its (now shifted)
# neighbour in column E so date cells keep the date format and the
# financial rows keep the numeric format, then the new value is written.
$ws.Range("E7").Copy($ws.Range("D7"))
$ws.Range("D7").Value = 43465
$ws.Range("E8").Copy($ws.Range("D8"))
$ws.Range("D8").Value = 1689700
$ws.Range("E9").Copy($ws.Range("D9"))
$ws.Range("D9").Value = "NA"
$ws.Range("E10").Copy($ws.Range("D10"))
$ws.Range("D10").Value = "NA"
$ws.Range("E12").Copy($ws.Range("D12"))
$ws.Range("D12").Value = "NA"
$ws.Range("E13").Copy($ws.Range("D13"))
$ws.Range("D13").Value = 0
$ws.Range("E14").Copy($ws.Range("D14"))
$ws.Range("D14").Value = 0
$ws.Range("E15").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 0
$ws.Range("E17").Copy($ws.Range("D17"))
$ws.Range("D17").Value = 676900
$ws.Range("E18").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 1012700
$ws.Range("E20").Copy($ws.Range("D20"))
$ws.Range("D20").Value = -455100
$ws.Range("E21").Copy($ws.Range("D21"))
$ws.Range("D21").Value = 590000
$ws.Range("E22").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 0
$ws.Range("E23").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 557700
$ws.Range("E24").Copy($ws.Range("D24"))
$ws.Range("D24").Value = 135300
$ws.Range("E25").Copy($ws.Range("D25"))
$ws.Range("D25").Value = 0
$ws.Range("E26").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 422400
$ws.Range("E27").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 384700
$ws.Range("E28").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 0
$ws.Range("E29").Copy($ws.Range("D29"))
$ws.Range("D29").Value = "NA"
$ws.Range("E30").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 0
$ws.Range("E31").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 0
$ws.Range("E32").Copy($ws.Range("D32"))
$ws.Range("D32").Value = 455100
$ws.Range("E33").Copy($ws.Range("D33"))
$ws.Range("D33").Value = 384700
$ws.Range("E34").Copy($ws.Range("D34"))
$ws.Range("D34").Value = 0
$ws.Range("E35").Copy($ws.Range("D35"))
$ws.Range("D35").Value = 384700
$ws.Range("E38").Copy($ws.Range("D38"))
$ws.Range("D38").Value = 43465
$ws.Range("E41").Copy($ws.Range("D41"))
$ws.Range("D41").Value = 1475000
$ws.Range("E42").Copy($ws.Range("D42"))
$ws.Range("D42").Value = 675100
$ws.Range("E43").Copy($ws.Range("D43"))
$ws.Range("D43").Value = 0
$ws.Range("E44").Copy($ws.Range("D44"))
$ws.Range("D44").Value = 0
$ws.Range("E45").Copy($ws.Range("D45"))
$ws.Range("D45").Value = 0
$ws.Range("E46").Copy($ws.Range("D46"))
$ws.Range("D46").Value = 0
$ws.Range("E47").Copy($ws.Range("D47"))
$ws.Range("D47").Value = 0
$ws.Range("E48").Copy($ws.Range("D48"))
$ws.Range("D48").Value = 346200
$ws.Range("E49").Copy($ws.Range("D49"))
$ws.Range("D49").Value = 2436900
$ws.Range("E50").Copy($ws.Range("D50"))
$ws.Range("D50").Value = 0
$ws.Range("E51").Copy($ws.Range("D51"))
$ws.Range("D51").Value = 0
$ws.Range("E52").Copy($ws.Range("D52"))
$ws.Range("D52").Value = 0
$ws.Range("E53").Copy($ws.Range("D53"))
$ws.Range("D53").Value = 0
$ws.Range("E54").Copy($ws.Range("D54"))
$ws.Range("D54").Value = 51899400
$ws.Range("E57").Copy($ws.Range("D57"))
$ws.Range("D57").Value = 0
$ws.Range("E58").Copy($ws.Range("D58"))
$ws.Range("D58").Value = 0
$ws.Range("E59").Copy($ws.Range("D59"))
$ws.Range("D59").Value = 0
$ws.Range("E60").Copy($ws.Range("D60"))
$ws.Range("D60").Value = 0
$ws.Range("E61").Copy($ws.Range("D61"))
$ws.Range("D61").Value = 954200
$ws.Range("E62").Copy($ws.Range("D62"))
$ws.Range("D62").Value = 79100
$ws.Range("E63").Copy($ws.Range("D63"))
$ws.Range("D63").Value = 0
$ws.Range("E64").Copy($ws.Range("D64"))
$ws.Range("D64").Value = 0
$ws.Range("E65").Copy($ws.Range("D65"))
$ws.Range("D65").Value = 0
$ws.Range("E66").Copy($ws.Range("D66"))
$ws.Range("D66").Value = 45244100
$ws.Range("E68").Copy($ws.Range("D68"))
$ws.Range("D68").Value = 0
$ws.Range("E69").Copy($ws.Range("D69"))
$ws.Range("D69").Value = 0
$ws.Range("E70").Copy($ws.Range("D70"))
$ws.Range("D70").Value = 502800
$ws.Range("E71").Copy($ws.Range("D71"))
$ws.Range("D71").Value = 0
$ws.Range("E72").Copy($ws.Range("D72"))
$ws.Range("D72").Value = 297200
$ws.Range("E73").Copy($ws.Range("D73"))
$ws.Range("D73").Value = 0
$ws.Range("E74").Copy($ws.Range("D74"))
$ws.Range("D74").Value = 0
$ws.Range("E75").Copy($ws.Range("D75"))
$ws.Range("D75").Value = 0
$ws.Range("E76").Copy($ws.Range("D76"))
$ws.Range("D76").Value = 6152400
$ws.Range("E77").Copy($ws.Range("D77"))
$ws.Range("D77").Value = 0
$ws.Range("E80").Copy($ws.Range("D80"))
$ws.Range("D80").Value = 43465
$ws.Range("E81").Copy($ws.Range("D81"))
$ws.Range("D81").Value = 384700
$ws.Range("E83").Copy($ws.Range("D83"))
$ws.Range("D83").Value = 32300
$ws.Range("E84").Copy($ws.Range("D84"))
$ws.Range("D84").Value = 0
$ws.Range("E85").Copy($ws.Range("D85"))
$ws.Range("D85").Value = 0
$ws.Range("E86").Copy($ws.Range("D86"))
$ws.Range("D86").Value = 0
$ws.Range("E87").Copy($ws.Range("D87"))
$ws.Range("D87").Value = 0
$ws.Range("E88").Copy($ws.Range("D88"))
$ws.Range("D88").Value = 0
$ws.Range("E89").Copy($ws.Range("D89"))
$ws.Range("D89").Value = 540400
$ws.Range("E91").Copy($ws.Range("D91"))
$ws.Range("D91").Value = -9800
$ws.Range("E92").Copy($ws.Range("D92"))
$ws.Range("D92").Value = 0
$ws.Range("E93").Copy($ws.Range("D93"))
$ws.Range("D93").Value = 0
$ws.Range("E94").Copy($ws.Range("D94"))
$ws.Range("D94").Value = -4020500
$ws.Range("E96").Copy($ws.Range("D96"))
$ws.Range("D96").Value = -365900
$ws.Range("E97").Copy($ws.Range("D97"))
$ws.Range("D97").Value = 0
$ws.Range("E98").Copy($ws.Range("D98"))
$ws.Range("D98").Value = 0
$ws.Range("E99").Copy($ws.Range("D99"))
$ws.Range("D99").Value = 0
$ws.Range("E100").Copy($ws.Range("D100"))
$ws.Range("D100").Value = 2426900
$ws.Range("E101").Copy($ws.Range("D101"))
$ws.Range("D101").Value = 0
$ws.Range("E102").Copy($ws.Range("D102"))
$ws.Range("D102").Value = -1053200
